# Updates cryptos list price (D) and volume-change (E) columns to the latest scrape.
# Source data is stored as literal text (inlineStr) even when it parses as a number
# (e.g. "1.00", "601.42"), so any cell whose new value would be auto-coerced to a
# number by Excel is first forced to Text format, written, then has its number format
# reset back to General via a Paste Special (Formats-only) from an untouched, default-
# styled cell -- this keeps the text value without leaving a stray cell style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text does not parse as a plain number: safe to assign directly ---
$directValues = @{
    "D2" = "65.775.53"
    "E2" = "  +0.35%  "
    "D3" = "2.673.96"
    "E3" = "  +0.89%  "
    "E4" = "  -0.01%  "
    "E5" = "  -0.63%  "
    "E6" = "  +0.56%  "
    "E7" = "  +0.02%  "
    "E8" = "  +5.42%  "
    "E9" = "  +0.68%  "
    "E10" = "  +1.01%  "
    "E11" = "  -0.29%  "
    "E12" = "  -0.17%  "
    "E13" = "  -1.07%  "
    "E14" = "  +1.03%  "
    "D15" = "3.152.81"
    "E15" = "  +0.82%  "
    "D16" = "65.574.05"
    "E16" = "  +0.43%  "
    "D17" = "2.671.37"
    "E17" = "  +0.37%  "
    "E18" = "  -0.48%  "
    "E19" = "  -1.25%  "
    "E20" = "  +1.27%  "
    "E21" = "  -1.84%  "
    "E22" = "  -0.10%  "
    "E23" = "  +0.32%  "
    "E24" = "  +5.39%  "
    "E25" = "  +3.98%  "
    "E26" = "  -4.52%  "
    "E27" = "  +1.80%  "
    "E28" = "  -1.99%  "
    "E29" = "  +0.22%  "
    "E30" = "  +3.29%  "
    "E31" = "  +0.48%  "
    "E32" = "  +0.29%  "
    "E33" = "  +0.04%  "
    "E34" = "  +3.92%  "
    "E35" = "  -1.13%  "
    "E36" = "  -1.65%  "
    "E37" = "  -1.21%  "
    "E38" = "  -0.01%  "
    "E39" = "  -2.49%  "
    "E40" = "  -1.15%  "
    "E41" = "  +0.01%  "
    "E42" = "  +1.93%  "
    "E43" = "  -0.13%  "
    "E44" = "  -0.83%  "
    "E45" = "  +1.30%  "
    "E46" = "  -1.56%  "
    "E47" = "  +1.67%  "
    "E48" = "  -1.00%  "
    "E49" = "  -1.15%  "
    "E50" = "  +3.41%  "
    "E51" = "  +3.13%  "
}
foreach ($ref in $directValues.Keys) {
    $ws.Range($ref).Value = $directValues[$ref]
}

# --- Cells whose new text WOULD be auto-converted to a number (e.g. "1.00" -> 1): ---
# --- force Text format, assign, then restore the default style.                 ---
$textValues = @{
    "D5" = "601.42"
    "D6" = "157.24"
    "D8" = "0.620"
    "D9" = "0.124"
    "D10" = "5.93"
    "D13" = "29.52"
    "D14" = "0.0000198"
    "D18" = "12.66"
    "D21" = "352.22"
    "D23" = "69.92"
    "D27" = "0.168"
    "D28" = "1.60"
    "D30" = "545.05"
    "D31" = "1.00"
    "D35" = "5.46"
    "D38" = "1.00"
    "D39" = "158.30"
    "D42" = "42.71"
    "D43" = "165.14"
    "D44" = "4.09"
    "D45" = "0.0616"
    "D46" = "2.31"
    "D47" = "23.32"
    "D48" = "0.645"
    "D50" = "0.101"
}
foreach ($ref in $textValues.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $textValues[$ref]
}

# A1 is blank and carries the sheet default style; use it as the Paste Special source
# so the forced Text format above does not leave a residual cell style.
$ws.Range("A1").Copy()
foreach ($ref in $textValues.Keys) {
    $ws.Range($ref).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

